$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.000364705512765795
$ws.Range("A3").Value = 0.00020177781698293984
$ws.Range("H3").Value = 5.396225929260254
$ws.Range("A4").Value = 0.00016292768123093992
$ws.Range("H4").Value = 5.288824081420898
